$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 12.88963960803372
$ws.Range("C2").Value = 9.976687499281631
$ws.Range("D2").Value = 6.048845672422907
$ws.Range("E2").Value = 12.85754554660048
$ws.Range("F2").Value = 29.05216494537615
$ws.Range("K2").Value = 9.279679946302073
$ws.Range("L2").Value = 9.936603969692772
$ws.Range("M2").Value = 14.34080071758721
$ws.Range("N2").Value = 20.61222577438471
$ws.Range("O2").Value = 26.07061801874621
$ws.Range("B3").Value = 12.68735875500144
$ws.Range("C3").Value = 9.968956650452613
$ws.Range("D3").Value = 6.005592915366127
$ws.Range("E3").Value = 12.88524249874731
$ws.Range("F3").Value = 29.08328807432049
$ws.Range("K3").Value = 9.125311463440566
$ws.Range("L3").Value = 9.943880510330171
$ws.Range("M3").Value = 14.31423113161551
$ws.Range("N3").Value = 20.67265565931606
$ws.Range("O3").Value = 26.13503928559392
$ws.Range("B4").Value = 12.56411785047949
$ws.Range("C4").Value = 9.964333342000804
$ws.Range("D4").Value = 5.978386136169828
$ws.Range("E4").Value = 12.90391149391134
$ws.Range("F4").Value = 29.10890588684179
$ws.Range("K4").Value = 9.030839042329019
$ws.Range("L4").Value = 9.949666414332624
$ws.Range("M4").Value = 14.30005469857685
$ws.Range("N4").Value = 20.7115067768385
$ws.Range("O4").Value = 26.17937452268554
$ws.Range("B5").Value = 12.51420271716208
$ws.Range("C5").Value = 9.962480838347759
$ws.Range("D5").Value = 5.967139253591748
$ws.Range("E5").Value = 12.91193771650268
$ws.Range("F5").Value = 29.12098065955844
$ws.Range("K5").Value = 8.992468489631881
$ws.Range("L5").Value = 9.952356223714142
$ws.Range("M5").Value = 14.29481970384186
$ws.Range("N5").Value = 20.72777945805703
$ws.Range("O5").Value = 26.19864222747058
$ws.Range("B6").Value = 12.50593476250154
$ws.Range("C6").Value = 9.962175150122537
$ws.Range("D6").Value = 5.965262159746192
$ws.Range("E6").Value = 12.91329574662445
$ws.Range("F6").Value = 29.12308438754117
$ws.Range("K6").Value = 8.986106261555211
$ws.Range("L6").Value = 9.952822932144295
$ws.Range("M6").Value = 14.29398329926264
$ws.Range("N6").Value = 20.73050816862602
$ws.Range("O6").Value = 26.20191409592297
$ws.Range("B7").Value = 12.56344335042476
$ws.Range("C7").Value = 9.964308230174725
$ws.Range("D7").Value = 5.97823509969812
$ws.Range("E7").Value = 12.90401804366376
$ws.Range("F7").Value = 29.10906211221668
$ws.Range("K7").Value = 9.03032098058512
$ws.Range("L7").Value = 9.949701345059252
$ws.Range("M7").Value = 14.29998189726069
$ws.Range("N7").Value = 20.71172445053245
$ws.Range("O7").Value = 26.17962951379204
$ws.Range("B8").Value = 12.81973282656905
$ws.Range("C8").Value = 9.973996568404276
$ws.Range("D8").Value = 6.034068932748896
$ws.Range("E8").Value = 12.86675048919323
$ws.Range("F8").Value = 29.06154489997011
$ws.Range("K8").Value = 9.226418906103932
$ws.Range("L8").Value = 9.93883970998168
$ws.Range("M8").Value = 14.33119874132577
$ws.Range("N8").Value = 20.63270018503911
$ws.Range("O8").Value = 26.09183756027815
$ws.Range("B9").Value = 13.3270945496728
$ws.Range("C9").Value = 9.993954380967484
$ws.Range("D9").Value = 6.13826038634759
$ws.Range("E9").Value = 12.80685058397548
$ws.Range("F9").Value = 29.02003287431345
$ws.Range("K9").Value = 9.611278558716329
$ws.Range("L9").Value = 9.927971913374147
$ws.Range("M9").Value = 14.40916008754214
$ws.Range("N9").Value = 20.49153649116023
$ws.Range("O9").Value = 25.95766272524025
$ws.Range("B10").Value = 13.69902152920693
$ws.Range("C10").Value = 10.00917540615244
$ws.Range("D10").Value = 6.211389967608484
$ws.Range("E10").Value = 12.77085983818215
$ws.Range("F10").Value = 29.021039377737
$ws.Range("K10").Value = 9.891413959447711
$ws.Range("L10").Value = 9.926308670798612
$ws.Range("M10").Value = 14.47633927577535
$ws.Range("N10").Value = 20.3961547583261
$ws.Range("O10").Value = 25.8823112019524
$ws.Range("B11").Value = 13.86727499664404
$ws.Range("C11").Value = 10.01621531379808
$ws.Range("D11").Value = 6.243876626081969
$ws.Range("E11").Value = 12.7562235385617
$ws.Range("F11").Value = 29.02832593483103
$ws.Range("K11").Value = 10.01771776489023
$ws.Range("L11").Value = 9.926915129298822
$ws.Range("M11").Value = 14.50897808471901
$ws.Range("N11").Value = 20.35455442222795
$ws.Range("O11").Value = 25.85308767151032
$ws.Range("B12").Value = 13.93079220390267
$ws.Range("C12").Value = 10.01889727425264
$ws.Range("D12").Value = 6.256062959058224
$ws.Range("E12").Value = 12.75093043620079
$ws.Range("F12").Value = 29.03206479633539
$ws.Range("K12").Value = 10.06533809156038
$ws.Range("L12").Value = 9.92733984248283
$ws.Range("M12").Value = 14.52162974193787
$ws.Range("N12").Value = 20.33905744786917
$ws.Range("O12").Value = 25.84274903646164
$ws.Range("B13").Value = 13.91712234487892
$ws.Range("C13").Value = 10.01831895883504
$ws.Range("D13").Value = 6.253443612382474
$ws.Range("E13").Value = 12.75205931570897
$ws.Range("F13").Value = 29.03121604295854
$ws.Range("K13").Value = 10.05509216381075
$ws.Range("L13").Value = 9.927239713440995
$ws.Range("M13").Value = 14.51889210489461
$ws.Range("N13").Value = 20.34238362318643
$ws.Range("O13").Value = 25.84494327030225
$ws.Range("B14").Value = 13.87250488374646
$ws.Range("C14").Value = 10.01643564079572
$ws.Range("D14").Value = 6.244881540530302
$ws.Range("E14").Value = 12.75578307634694
$ws.Range("F14").Value = 29.02861391667388
$ws.Range("K14").Value = 10.02163994347132
$ws.Range("L14").Value = 9.926946167106424
$ws.Range("M14").Value = 14.51001313151979
$ws.Range("N14").Value = 20.35327434921377
$ws.Range("O14").Value = 25.8522225174665
$ws.Range("B15").Value = 13.84514798691633
$ws.Range("C15").Value = 10.01528413442727
$ws.Range("D15").Value = 6.239621860829484
$ws.Range("E15").Value = 12.75809645035616
$ws.Range("F15").Value = 29.02714752961196
$ws.Range("K15").Value = 10.00112106539752
$ws.Range("L15").Value = 9.926791734697114
$ws.Range("M15").Value = 14.50461232820576
$ws.Range("N15").Value = 20.35997856098112
$ws.Range("O15").Value = 25.85677605567253
$ws.Range("B16").Value = 13.68800155251715
$ws.Range("C16").Value = 10.00871761468968
$ws.Range("D16").Value = 6.209250889392368
$ws.Range("E16").Value = 12.77185126393949
$ws.Range("F16").Value = 29.02070041008113
$ws.Range("K16").Value = 9.88313300834214
$ws.Range("L16").Value = 9.926296386572323
$ws.Range("M16").Value = 14.47424750179616
$ws.Range("N16").Value = 20.39890934659038
$ws.Range("O16").Value = 25.88432281238601
$ws.Range("B17").Value = 13.59131230068103
$ws.Range("C17").Value = 10.00471853392096
$ws.Range("D17").Value = 6.190416943551803
$ws.Range("E17").Value = 12.78073384018858
$ws.Range("F17").Value = 29.0184928584344
$ws.Range("K17").Value = 9.810428399996345
$ws.Range("L17").Value = 9.926340972557551
$ws.Range("M17").Value = 14.45614728170185
$ws.Range("N17").Value = 20.42324959729114
$ws.Range("O17").Value = 25.90251705332294
$ws.Range("B18").Value = 13.53561446832655
$ws.Range("C18").Value = 10.00242930009708
$ws.Range("D18").Value = 6.179510907631435
$ws.Range("E18").Value = 12.78600628945225
$ws.Range("F18").Value = 29.01786606940164
$ws.Range("K18").Value = 9.768506757408129
$ws.Range("L18").Value = 9.926494977936068
$ws.Range("M18").Value = 14.44593259868132
$ws.Range("N18").Value = 20.43741794912275
$ws.Range("O18").Value = 25.9134575870595
$ws.Range("B19").Value = 13.51674351047339
$ws.Range("C19").Value = 10.00165609904
$ws.Range("D19").Value = 6.175805836776846
$ws.Range("E19").Value = 12.78781952756699
$ws.Range("F19").Value = 29.01776433418555
$ws.Range("K19").Value = 9.754296389315057
$ws.Range("L19").Value = 9.926569195225307
$ws.Range("M19").Value = 14.44250796396719
$ws.Range("N19").Value = 20.44224408135466
$ws.Range("O19").Value = 25.91724353396012
$ws.Range("B20").Value = 13.60161426178778
$ws.Range("C20").Value = 10.00514311547727
$ws.Range("D20").Value = 6.192429455323825
$ws.Range("E20").Value = 12.77977136287031
$ws.Range("F20").Value = 29.01866132932367
$ws.Range("K20").Value = 9.818179016577126
$ws.Range("L20").Value = 9.926322948081367
$ws.Range("M20").Value = 14.4580538356476
$ws.Range("N20").Value = 20.42064110684913
$ws.Range("O20").Value = 25.90053100649056
$ws.Range("B21").Value = 13.88561593346206
$ws.Range("C21").Value = 10.01698838441084
$ws.Range("D21").Value = 6.247399594691681
$ws.Range("E21").Value = 12.75468255286407
$ws.Range("F21").Value = 29.02935166058636
$ws.Range("K21").Value = 10.03147167397614
$ws.Range("L21").Value = 9.927027102828799
$ws.Range("M21").Value = 14.51261323137925
$ws.Range("N21").Value = 20.35006853553772
$ws.Range("O21").Value = 25.85006467080281
$ws.Range("B22").Value = 14.07005385657072
$ws.Range("C22").Value = 10.02482353779251
$ws.Range("D22").Value = 6.282650053612063
$ws.Range("E22").Value = 12.73973879542503
$ws.Range("F22").Value = 29.04204665919802
$ws.Range("K22").Value = 10.16963768517806
$ws.Range("L22").Value = 9.928623776861663
$ws.Range("M22").Value = 14.54996982768572
$ws.Range("N22").Value = 20.30543792004244
$ws.Range("O22").Value = 25.82132378295466
$ws.Range("B23").Value = 13.97174303331895
$ws.Range("C23").Value = 10.0206333889794
$ws.Range("D23").Value = 6.263899169006613
$ws.Range("E23").Value = 12.74758169317035
$ws.Range("F23").Value = 29.0347497628095
$ws.Range("K23").Value = 10.09602320546048
$ws.Range("L23").Value = 9.927667947914095
$ws.Range("M23").Value = 14.52987881043382
$ws.Range("N23").Value = 20.32912191275486
$ws.Range("O23").Value = 25.83627496251946
$ws.Range("B24").Value = 13.59695708272407
$ws.Range("C24").Value = 10.00495113131033
$ws.Range("D24").Value = 6.191519841906291
$ws.Range("E24").Value = 12.78020598259297
$ws.Range("F24").Value = 29.01858316255045
$ws.Range("K24").Value = 9.814675341808032
$ws.Range("L24").Value = 9.926330697013231
$ws.Range("M24").Value = 14.45719128598016
$ws.Range("N24").Value = 20.42181986073427
$ws.Range("O24").Value = 25.90142740183839
$ws.Range("B25").Value = 13.18971984047837
$ws.Range("C25").Value = 9.988455321295408
$ws.Range("D25").Value = 6.110660836300553
$ws.Range("E25").Value = 12.82164571916699
$ws.Range("F25").Value = 29.02572806820231
$ws.Range("K25").Value = 9.507425693058201
$ws.Range("L25").Value = 9.92979879955856
$ws.Range("M25").Value = 14.3863070898689
$ws.Range("N25").Value = 20.52825598553969
$ws.Range("O25").Value = 25.989886235806
